$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 - "Khoi luong" values
$ws.Range("A3").Value = 0
$ws.Range("B3").Value = 655
$ws.Range("C3").Value = 1333
$ws.Range("D3").Value = 9096

# Row 4 - "Khoi luong" values
$ws.Range("A4").Value = 121.85
$ws.Range("B4").Value = 1291.6500000000001
$ws.Range("C4").Value = 4687.95
$ws.Range("D4").Value = 14931.5

# Row 6 - "RSI" values
$ws.Range("A6").Value = 40.25
$ws.Range("B6").Value = 31.77
$ws.Range("C6").Value = 29.23
$ws.Range("D6").Value = 26.66

# Row 7 - "RSI" values
$ws.Range("A7").Value = 41.37
$ws.Range("B7").Value = 32.28
$ws.Range("C7").Value = 35.57
$ws.Range("D7").Value = 34.56

$excel.CalculateFull()
